$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows: correct a typo and fill in missing DK P/N / Supplier cells ---

# Row 9 (LX): add DigiKey part number
$ws.Range("F9").Value = "490-14207-1-ND"

# Row 10 (U1): fix part number typo and add DigiKey part number
$ws.Range("E10").Value = "BD9G101G"
$ws.Range("F10").Value = "BD9G101G-CT-ND"

# Row 11 (CIN2): add Supplier and DigiKey part number
$ws.Range("D11").Value = "ROHM"
$ws.Range("F11").Value = "SMAJ30ALFCT-ND"

# --- Insert three new BOM rows (JEN, RFRA, RFB) after row 11 ---
$ws.Range("A12:A14").EntireRow.Insert()

# Copy formatting from the row above (CIN2, row 11) onto the new rows
$ws.Range("A11:G11").Copy()
$ws.Range("A12:G14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A12").Value = "JEN"
$ws.Range("B12").Value = "short"
$ws.Range("C12").Value = "'0603"
$ws.Range("D12").Value = "ROHM"
$ws.Range("E12").Value = "MCR03ERTJ000"
$ws.Range("F12").Value = "RHM0.0CGTR-ND"

$ws.Range("A13").Value = "RFRA"
$ws.Range("B13").Value = "short"
$ws.Range("C13").Value = "'0603"
$ws.Range("D13").Value = "ROHM"
$ws.Range("E13").Value = "MCR03ERTJ000"
$ws.Range("F13").Value = "RHM0.0CGTR-ND"

$ws.Range("A14").Value = "RFB"
$ws.Range("B14").Value = "short"
$ws.Range("C14").Value = "'0603"
$ws.Range("D14").Value = "ROHM"
$ws.Range("E14").Value = "MCR03ERTJ000"
$ws.Range("F14").Value = "RHM0.0CGTR-ND"

# --- Fill in the Test Pin rows (now rows 15 & 16 after the insert) ---

# Row 15: Vcc, Vout test pin
$ws.Range("B15").Value = "TP"
$ws.Range("C15").Value = "Via"
$ws.Range("D15").Value = "Keystone"
$ws.Range("F15").Value = "36-5005-ND"

# Row 16: GND test pin
$ws.Range("B16").Value = "TP"
$ws.Range("C16").Value = "Via"
$ws.Range("D16").Value = "Keystone"
$ws.Range("F16").Value = "36-5006-ND"

# --- Restore the active selection ---
$ws.Range("K11").Select()
